# Tutorial 6 solution update:
# - Reformat the attendance dates in column A from DD/MM/YYYY to DD-MM-YYYY
# - Recompute the Total Attendance Count / Real / Duplicate / Invalid / Absent
#   columns (D-H) for the rows whose tallies changed as part of the update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DateText($row, $text) {
    # Force the cell to be treated as plain text so Excel does not
    # auto-convert the DD-MM-YYYY looking string into a date serial
    # number, then restore the cell's original (default) formatting.
    $cell = $ws.Range("A$row")
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

Set-DateText 3  "28-07-2022"
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1

Set-DateText 4  "01-08-2022"
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("H4").Value = 0

Set-DateText 5  "04-08-2022"
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("H5").Value = 0

Set-DateText 6  "08-08-2022"
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1
$ws.Range("H6").Value = 0

Set-DateText 7  "11-08-2022"

Set-DateText 8  "15-08-2022"

Set-DateText 9  "18-08-2022"

Set-DateText 10 "22-08-2022"

Set-DateText 11 "25-08-2022"

Set-DateText 12 "29-08-2022"
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 1
$ws.Range("H12").Value = 0

Set-DateText 13 "01-09-2022"

Set-DateText 14 "05-09-2022"
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 1
$ws.Range("H14").Value = 0

Set-DateText 15 "08-09-2022"
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 1
$ws.Range("H15").Value = 0

Set-DateText 16 "12-09-2022"

Set-DateText 17 "15-09-2022"

Set-DateText 18 "19-09-2022"

Set-DateText 19 "22-09-2022"

Set-DateText 20 "26-09-2022"

Set-DateText 21 "29-09-2022"
